$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Overview")
$ws.Range("B2").Value = "4/28/2023, 1:46:02 PM CDT"

$ws = $wb.Worksheets.Item("Defs")
$ws.Range("A2").Value = "lh0wmwkt-vnhp"
$ws.Range("B2").Value = "4/28/2023, 1:46:02 PM"
$ws.Range("C2").Value = "4/28/2023, 1:46:02 PM"
$ws.Range("A3").Value = "lh0wmwkx-0sbz"
$ws.Range("B3").Value = "4/28/2023, 1:46:02 PM"
$ws.Range("C3").Value = "4/28/2023, 1:46:02 PM"
$ws.Range("A4").Value = "lh0wmwky-q0sf"
$ws.Range("C4").Value = "4/28/2023, 1:46:02 PM"
$ws.Range("A5").Value = "lh0wmwlk-xg1i"
$ws.Range("B5").Value = "4/28/2023, 1:46:02 PM"
$ws.Range("C5").Value = "4/28/2023, 1:46:02 PM"

$ws = $wb.Worksheets.Item("Point Defs")
$ws.Range("A2").Value = "lh0wmwkw-xm25"
$ws.Range("B2").Value = "4/28/2023, 1:46:02 PM"
$ws.Range("C2").Value = "4/28/2023, 1:46:02 PM"
$ws.Range("A3").Value = "lh0wmwky-ckzr"
$ws.Range("B3").Value = "4/28/2023, 1:46:02 PM"
$ws.Range("C3").Value = "4/28/2023, 1:46:02 PM"
$ws.Range("F3").Value = "0f3n"
$ws.Range("A4").Value = "lh0wmwkz-heaj"
$ws.Range("B4").Value = "4/28/2023, 1:46:02 PM"
$ws.Range("C4").Value = "4/28/2023, 1:46:02 PM"
$ws.Range("J4").Value = "NUMBER"
$ws.Range("A5").Value = "lh0wmwkz-iqss"
$ws.Range("B5").Value = "4/28/2023, 1:46:02 PM"
$ws.Range("C5").Value = "4/28/2023, 1:46:02 PM"
$ws.Range("A6").Value = "lh0wmwll-7w4j"
$ws.Range("B6").Value = "4/28/2023, 1:46:02 PM"
$ws.Range("C6").Value = "4/28/2023, 1:46:02 PM"
$ws.Range("J6").Value = "NUMBER"
$ws.Range("A7").Value = "lh0wmwll-pg69"
$ws.Range("B7").Value = "4/28/2023, 1:46:02 PM"
$ws.Range("C7").Value = "4/28/2023, 1:46:02 PM"

$ws = $wb.Worksheets.Item("Entry")
$ws.Range("A2").Value = "lh0wmwl0-zqz7"
$ws.Range("B2").Value = "4/28/2023, 1:46:02 PM"
$ws.Range("C2").Value = "4/28/2023, 1:46:02 PM"
$ws.Range("A3").Value = "lh0wmwl1-y0d3"
$ws.Range("B3").Value = "4/28/2023, 1:46:02 PM"
$ws.Range("C3").Value = "4/28/2023, 1:46:02 PM"
$ws.Range("F3").Value = "lh0wmwl2-bkai"
$ws.Range("G3").Value = "2023-04-28T13:46:02"
$ws.Range("A4").Value = "lh0wmwlm-sulo"
$ws.Range("B4").Value = "4/28/2023, 1:46:02 PM"
$ws.Range("C4").Value = "4/28/2023, 1:46:02 PM"

$ws = $wb.Worksheets.Item("Entry Points")
$ws.Range("A2").Value = "lh0wmwl1-crrp"
$ws.Range("B2").Value = "4/28/2023, 1:46:02 PM"
$ws.Range("C2").Value = "4/28/2023, 1:46:02 PM"
$ws.Range("A3").Value = "lh0wmwl1-07ev"
$ws.Range("B3").Value = "4/28/2023, 1:46:02 PM"
$ws.Range("C3").Value = "4/28/2023, 1:46:02 PM"
$ws.Range("A4").Value = "lh0wmwlm-84du"
$ws.Range("B4").Value = "4/28/2023, 1:46:02 PM"
$ws.Range("C4").Value = "4/28/2023, 1:46:02 PM"
$ws.Range("A5").Value = "lh0wmwln-580r"
$ws.Range("B5").Value = "4/28/2023, 1:46:02 PM"
$ws.Range("C5").Value = "4/28/2023, 1:46:02 PM"

$ws = $wb.Worksheets.Item("Tag Defs")
$ws.Range("A2").Value = "lh0wmwl2-036u"
$ws.Range("B2").Value = "4/28/2023, 1:46:02 PM"
$ws.Range("C2").Value = "4/28/2023, 1:46:02 PM"
$ws.Range("E2").Value = "plp7"
$ws.Range("A3").Value = "lh0wmwl3-6xxo"
$ws.Range("B3").Value = "4/28/2023, 1:46:02 PM"
$ws.Range("C3").Value = "4/28/2023, 1:46:02 PM"
$ws.Range("A4").Value = "lh0wmwl4-0bew"
$ws.Range("B4").Value = "4/28/2023, 1:46:02 PM"
$ws.Range("C4").Value = "4/28/2023, 1:46:02 PM"
$ws.Range("A5").Value = "lh0wmwln-0qmh"
$ws.Range("B5").Value = "4/28/2023, 1:46:02 PM"
$ws.Range("C5").Value = "4/28/2023, 1:46:02 PM"

$ws = $wb.Worksheets.Item("Tags")
$ws.Range("A2").Value = "lh0wmwl4-1lqh"
$ws.Range("B2").Value = "4/28/2023, 1:46:02 PM"
$ws.Range("C2").Value = "4/28/2023, 1:46:02 PM"
$ws.Range("A3").Value = "lh0wmwl5-0l5y"
$ws.Range("B3").Value = "4/28/2023, 1:46:02 PM"
$ws.Range("C3").Value = "4/28/2023, 1:46:02 PM"
